$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 620, pushing existing rows 620-679 down to 621-680.
$ws.Rows.Item(620).Insert()

# Populate the new row 620 with a fresh record (same constant columns as the
# rest of this "Pomelo" block, with a new date / quality / volume / prices).
$ws.Cells.Item(620, 1).Value = 4
$ws.Cells.Item(620, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(620, 3).Value = "Los Lagos"
$ws.Cells.Item(620, 4).Value = 45194
$ws.Cells.Item(620, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(620, 5).Value = 10
$ws.Cells.Item(620, 6).Value = "Fruta"
$ws.Cells.Item(620, 7).Value = 100102
$ws.Cells.Item(620, 8).Value = "Cítricos"
$ws.Cells.Item(620, 9).Value = 100102006
$ws.Cells.Item(620, 10).Value = "Pomelo"
$ws.Cells.Item(620, 11).Value = "Start Ruby"
$ws.Cells.Item(620, 12).Value = "Primera"
$ws.Cells.Item(620, 13).Value = 100
$ws.Cells.Item(620, 14).Value = 15000
$ws.Cells.Item(620, 15).Value = 16000
$ws.Cells.Item(620, 16).Value = 15500
$ws.Cells.Item(620, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(620, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(620, 19).Value = 1107
$ws.Cells.Item(620, 20).Value = 14
